$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 152 (pushing old rows 152-166 down to 155-169).
# This mirrors a weekly data refresh: three new observations are prepended
# to this producer's "Papa" price series.
$ws.Rows.Item(152).Resize(3).Insert()

# Row 152: new weekly price entry
$ws.Range("A152").Value = 1
$ws.Range("B152").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C152").Value = "Arica y Parinacota"
$ws.Range("D152").Value = 44984
$ws.Range("E152").Value = 15
$ws.Range("F152").Value = 100114001
$ws.Range("G152").Value = "Papa"
$ws.Range("H152").Value = "Asterix"
$ws.Range("I152").Value = "1a (cosecha lavada)"
$ws.Range("J152").Value = 900
$ws.Range("K152").Value = 16000
$ws.Range("L152").Value = 17000
$ws.Range("M152").Value = 16611
$ws.Range("N152").Value = "`$/malla 25 kilos"
$ws.Range("O152").Value = "Región de Los Lagos"
$ws.Range("P152").Value = 664
$ws.Range("Q152").Value = 25
$ws.Range("R152").Value = "Hortaliza"

# Row 153: new weekly price entry
$ws.Range("A153").Value = 1
$ws.Range("B153").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C153").Value = "Arica y Parinacota"
$ws.Range("D153").Value = 44984
$ws.Range("E153").Value = 15
$ws.Range("F153").Value = 100114001
$ws.Range("G153").Value = "Papa"
$ws.Range("H153").Value = "Asterix"
$ws.Range("I153").Value = "1a (cosecha)"
$ws.Range("J153").Value = 1000
$ws.Range("K153").Value = 14000
$ws.Range("L153").Value = 15000
$ws.Range("M153").Value = 14450
$ws.Range("N153").Value = "`$/saco 25 kilos"
$ws.Range("O153").Value = "Región de Los Lagos"
$ws.Range("P153").Value = 578
$ws.Range("Q153").Value = 25
$ws.Range("R153").Value = "Hortaliza"

# Row 154: new weekly price entry
$ws.Range("A154").Value = 1
$ws.Range("B154").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C154").Value = "Arica y Parinacota"
$ws.Range("D154").Value = 44984
$ws.Range("E154").Value = 15
$ws.Range("F154").Value = 100114001
$ws.Range("G154").Value = "Papa"
$ws.Range("H154").Value = "Red Lady"
$ws.Range("I154").Value = "1a (cosecha)"
$ws.Range("J154").Value = 950
$ws.Range("K154").Value = 13000
$ws.Range("L154").Value = 14000
$ws.Range("M154").Value = 13421
$ws.Range("N154").Value = "`$/saco 25 kilos"
$ws.Range("O154").Value = "Región del Bíobío"
$ws.Range("P154").Value = 537
$ws.Range("Q154").Value = 25
$ws.Range("R154").Value = "Hortaliza"
